# Romanian translation pass for CrisisText Video Scripts (Onboarding + Mindfulness)
$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 0, $false, $replace, 2) | Out-Null
}

# --- Onboarding lesson intro paragraph (two runs in the same <w:p>) ---
Replace-All "Welcome. This chatbot is here to help you prepare and feel more confident if things around you ever change or become uncertain. " "Bine ai venit. Busola Părintelui este aici pentru a te ajuta să fi pregătit și să te simți mai încrezător dacă vreodată lucrurile din jurul tău se vor schimba sau vor deveni incerte. "

Replace-All "The tips you will receive by using Busola Parintelui are meant to help you deal with uncertain or crisis times. We hope we can help you, even just a little bit, to be prepared " "Sfaturile pe care le vei primi prin intermediul chatbot-ului Busola Părintelui sunt menite să te ajute să faci față perioadelor de incertitudine sau criză. Sperăm că te putem ajuta, chiar și puțin, să te simți mai pregătit. "

Replace-All "The tips you receive here are designed with support from the World Vision Romania, World Health Organisation, UNICEF, UNHCR, and experts in the parenting research community. " "Sfaturile pe care le primiți aici sunt concepute cu sprijinul World Vision România, Organizației Mondiale a Sănătății, UNICEF, UNHCR și experților din comunitatea de cercetare în domeniul creșterii copiilor. "

# --- Standalone "Welcome" on-slide text cell: targeted via Paragraphs so the
#     other three "Welcome..." occurrences elsewhere are left untouched.
#     (Paragraph.Range.Text carries trailing paragraph/cell-mark control
#     characters, so trim those off before comparing.) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Welcome") {
        $para.Range.Text = "Bine ai venit"
        break
    }
}

Replace-All "Before you receive the tips, I want to explain how Busola Părintelui works. " "Înainte să primești sfaturile, vreau să îți explic cum funcționează Busola Părintelui. "

Replace-All "Together we will review: " "Împreună vom trece în revistă: "

# Longer variants containing "And " must be replaced before the shorter ones.
Replace-All "And how to access additional information & support" "Cum să accesezi informații suplimentare și asistență în situații de urgență"
Replace-All "How to access additional information & support" "Cum să accesezi informații suplimentare și asistență în situații de urgență"

# "What to expect" / "How to access playful activities..." each appear twice
# (identically) and both instances map to the same Romanian text.
Replace-All "What to expect" "La ce să te aștepți"
Replace-All "How to access playful activities with your child" "Cum să accesezi activități pentru tine și copilul tău"

Replace-All "You will receive a new tip every day to help you support your children in a time of crisis or uncertainty. " "Vei primi în fiecare zi un sfat nou care te va ajuta să-ți susții copiii în perioade de criză sau incertitudine. "

Replace-All "These tips are made to be short but helpful. It takes less than 5 minutes to read the daily tips. " "Aceste sfaturi sunt scurte, dar utile. Îți va lua mai puțin de 5 minute să citești sfaturile zilnice. "

# --- "Get more help" panel: replace full sentences before their short
#     standalone menu-item counterparts so substrings don't collide. ---
Replace-All "If you want to review any of the tips you’ve previously received, just type MENU and navigate to “Review Tips”" "Dacă vrei să recitești oricare dintre sfaturile primite anterior, tastează MENIU și accesează „Recitește Sfaturi”"

Replace-All "To change your language or gender settings, select “Change my Settings”" "Pentru a schimba setările de limbă sau gen, selectează „Modifică setările mele”"

Replace-All "To share this chatbot with a friend, select “Invite a Friend to Busola Părintelui”" "Pentru a recomanda acest chatbot unui prieten, selectează „Invită un prieten pe Busola Părintelui”"

Replace-All "For more information or resources available to you in times of crisis, select “Get more help.” You can also access this information by typing HELP at any time. " "Pentru mai multe informații sau resurse disponibile în situații de criză, selectează „Obține mai mult ajutor”. De asemenea, poți accesa aceste informații tastând cuvântul AJUTOR în orice moment. "

Replace-All "Finally, selecting “Watch a video about Busola Părintelui” will replay this video. " "În cele din urmă, selectând „Vizionați un videoclip despre Busola Părintelui” veți revedea acest videoclip. "

# --- Standalone MENU block (on-slide text) ---
Replace-All "MENU " "MENIU "
Replace-All "“What would you like to do?” " "„Ce ai vrei să faci?” "
Replace-All "Review Tips" "Recitește Sfaturi"
Replace-All "Change my Settings" "Modifică Setările Mele"
Replace-All "Invite a Friend to Busola Părintelui" "Invită un prieten pe Busola Părintelui"
Replace-All "Get more help" "Obține mai mult ajutor"
Replace-All "Watch a video about Busola Părintelui" "Vizionează un video despre Busola Părintelui"
Replace-All "Exit Menu" "Meniu de Ieșire"

# --- Mindfulness / playful activity tips ---
Replace-All "In a challenging time, it can be difficult to find moments to connect with our children, but these moments, even if they are small, give our children much-needed stability. " "Într-o perioadă dificilă, poate fi greu să găsim momente în care să ne conectăm cu copiii noștri, dar aceste momente, chiar dacă sunt scurte, le oferă copiilor noștri stabilitatea de care au atâta nevoie. "

Replace-All "Busola Părintelui offers ideas on how to play with your child. You can do these activities anywhere, without supplies. " "Busola Părintelui oferă idei despre cum să te joci cu copilul tău. Poți face aceste activități oriunde, fără să ai nevoie de materiale speciale. "

Replace-All "After each tip, you’ll be asked whether you’d like a playful activity, or to finish your lesson for the day. " "După fiecare sfat, vei fi întrebat dacă dorești să faci o activitate ludică sau să închei lecția pentru ziua respectivă. "

Replace-All "You can choose the type of activity you want to do with your children:" "Poți alege tipul de activitate pe care vrei să o faci alături de copiii tăi:"

Replace-All "Active - for energetic fun" "Activ - pentru distracție energică"
Replace-All "Calm - to relax together, or" "Liniștit - pentru a vă relaxa împreună sau"
Replace-All "Quick - for when you are short on time" "Rapid - pentru momentele în care aveți puțin timp"

Write-Host "done"
